# Commercial sector demand update and documentation update
#
# Adds a new "Documentation" sheet at the front of the workbook explaining
# the purpose of the "Uni" and "Bi" trade-link sheets, and removes the
# unused, empty "Sheet2".

$wb = $excel.ActiveWorkbook

# Insert the new Documentation sheet before the current first sheet ("Uni")
# so the final tab order is Documentation, Uni, Bi.
$doc = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$doc.Name = "Documentation"

$doc.Range("A1").Value = "Workbook: Defines all unilateral and/or bilateral trade links between regions"
$doc.Range("A2").Value = "Uni: Unilateral trade links"
$doc.Range("A3").Value = "Bi: Bilateral trade links"

# Remove the unused, empty "Sheet2"
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet2").Delete()
$excel.DisplayAlerts = $true

# Make the new Documentation sheet the active/selected sheet
$doc.Activate()
$doc.Range("A18").Select()
